$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 on the "Input" sheet is a stray blank divider row with no real
# data in it ("bad data"). Delete the whole row so every row below it
# (29-49, which already contains the real CO# log entries) shifts up by
# one - row 49's content becomes row 48, etc.
$ws.Rows.Item(28).Delete()

# Leave the selection where the saved workbook shows it.
$ws.Range("A4:XFD4").Select()
